# Invoice template: switch the numeric placeholder tokens (Qty, Unit Price,
# Tax Rate %) over to the "{=...}" computed/number syntax used by the
# templating engine, matching the other "{#...}" control tokens already in
# the sheet. Plain text placeholders stay as "{...}".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sales detail line: Qty (B16) and Unit Price (D16) placeholders.
$ws.Range("B16").Value = "{=qty}"
$ws.Range("D16").Value = "{=price}"

# Sales Tax row: Tax Rate % placeholder.
$ws.Range("E18").Value = "{=taxRatePercent}"
